# Populated all files with crossref tags, when applicable
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$authorsRow2 = "[Shuang-Jiang%Zhou%NULL%0,                   Li-Gang%Zhang%NULL%0,                   Lei-Lei%Wang%NULL%0,                   Zhao-Chang%Guo%NULL%0,                   Jing-Qi%Wang%NULL%0,                   Jin-Cheng%Chen%NULL%0,                   Mei%Liu%NULL%0,                   Xi%Chen%NULL%0,                   Jing-Xu%Chen%chenjx1110@163.com%0]"

$authorsRow3 = "[Benjamin%Oosterhoff%Benjamin.oosterhoff@montana.edu%0,                   Cara A.%Palmer%NULL%0,                   Jenna%Wilson%NULL%0,                   Natalie%Shook%NULL%0]"

$authorsRow4 = "[İsmail%Seçer%ismailsecer84@gmail.com%0,                   Sümeyye%Ulaş%NULL%0,                   Sümeyye%Ulaş%NULL%0]"

$authorsRow5 = "[Miao%Qu%NULL%1,                   Kun%Yang%NULL%1,                   Yujia%Cao%NULL%1,                   Mei Hong%Xiu%xiumeihong97@163.com%1,                   Xiang Yang%Zhang%zhangxy@psych.ac.cn%2,                   Xiang Yang%Zhang%zhangxy@psych.ac.cn%0]"

$crossrefSpringer = "_PMC_Springer_CROSSREF"
$crossrefElsevier = "_PMC_elsevier_CROSSREF"

# Column E = Authors, Column I = Other found locations
$ws.Range("E2").Value = $authorsRow2
$ws.Range("I2").Value = $crossrefSpringer

$ws.Range("E3").Value = $authorsRow3
$ws.Range("I3").Value = $crossrefElsevier

$ws.Range("E4").Value = $authorsRow4
$ws.Range("I4").Value = $crossrefSpringer

$ws.Range("E5").Value = $authorsRow5
$ws.Range("I5").Value = $crossrefSpringer
